$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns store plain text values
# (e.g. "27.555.50", "  +3.75%  ") rather than numbers. Excel's normal
# Value assignment auto-detects numeric-looking text and converts it to
# a real number, so we temporarily force the target range to Text format,
# write the literal strings, then restore the original General format/style
# so no stray formatting changes are left behind.
$targetRange = $ws.Range("D2:E51")
$targetRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.555.50'
$ws.Range("E2").Value = '  +3.75%  '
$ws.Range("D3").Value = '1.824.28'
$ws.Range("E3").Value = '  +4.69%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '342.78'
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = '0.3821'
$ws.Range("E7").Value = '  +1.17%  '
$ws.Range("D8").Value = '0.3539'
$ws.Range("E8").Value = '  +4.52%  '
$ws.Range("D9").Value = '49.91'
$ws.Range("E9").Value = '  +2.97%  '
$ws.Range("D10").Value = '1.238'
$ws.Range("E10").Value = '  +4.38%  '
$ws.Range("D11").Value = '0.07748'
$ws.Range("E11").Value = '  +3.61%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '22.30'
$ws.Range("E13").Value = '  +9.34%  '
$ws.Range("D14").Value = '6.625'
$ws.Range("E14").Value = '  +2.64%  '
$ws.Range("D15").Value = '1.825.59'
$ws.Range("E15").Value = '  +4.87%  '
$ws.Range("D16").Value = '7.238'
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("E17").Value = '  +3.68%  '
$ws.Range("D18").Value = '0.06748'
$ws.Range("D19").Value = '87.08'
$ws.Range("E19").Value = '  +4.40%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '17.61'
$ws.Range("E21").Value = '  +5.18%  '
$ws.Range("D22").Value = '6.540'
$ws.Range("E22").Value = '  +5.45%  '
$ws.Range("D23").Value = '13.17'
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("D24").Value = '27.541.20'
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("D25").Value = '2.483'
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("D26").Value = '2.684'
$ws.Range("E26").Value = '  +9.07%  '
$ws.Range("D27").Value = '22.06'
$ws.Range("E27").Value = '  +12.14%  '
$ws.Range("E28").Value = '  +5.60%  '
$ws.Range("D29").Value = '152.88'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("D30").Value = '2.029.84'
$ws.Range("D31").Value = '135.54'
$ws.Range("E31").Value = '  +2.73%  '
$ws.Range("D32").Value = '6.341'
$ws.Range("E32").Value = '  +3.86%  '
$ws.Range("D33").Value = '4.090'
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("D34").Value = '13.95'
$ws.Range("E34").Value = '  +7.43%  '
$ws.Range("D35").Value = '0.08798'
$ws.Range("E35").Value = '  +1.53%  '
$ws.Range("D36").Value = '1.698'
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = '5.638'
$ws.Range("E37").Value = '  +3.75%  '
$ws.Range("D38").Value = '0.7025'
$ws.Range("E38").Value = '  +12.64%  '
$ws.Range("D39").Value = '9.138'
$ws.Range("E39").Value = '  +6.23%  '
$ws.Range("D40").Value = '0.06525'
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("D41").Value = '0.2261'
$ws.Range("E41").Value = '  +3.56%  '
$ws.Range("D42").Value = '0.02408'
$ws.Range("E42").Value = '  +2.13%  '
$ws.Range("D43").Value = '1.306'
$ws.Range("E43").Value = '  +6.17%  '
$ws.Range("D44").Value = '14.81'
$ws.Range("E44").Value = '  +3.42%  '
$ws.Range("D45").Value = '0.6628'
$ws.Range("E45").Value = '  +9.17%  '
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Value = '3.944'
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("D48").Value = '2.193'
$ws.Range("E48").Value = '  +6.21%  '
$ws.Range("E49").Value = '  +3.68%  '
$ws.Range("D50").Value = '0.07312'
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").Value = '81.19'
$ws.Range("E51").Value = '  +4.26%  '

$targetRange.NumberFormat = "General"
$targetRange.Style = "Normal"

